$wb = $excel.ActiveWorkbook

# Row 51 (A Bile Business | Shark Oil)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 6719.107
$ws.Range("I51").Value = 10425.143
$ws.Range("J51").Value = 5483.7617
$ws.Range("K51").Value = 10425.143
$ws.Range("L51").Value = 5483.7617
$ws.Range("M51").Value = -9941.143
$ws.Range("N51").Value = -6451.7617

# Row 106 (Making Your Mark | Enchanted Palladium Ink)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 8096.9287
$ws.Range("I106").Value = 7504.385
$ws.Range("J106").Value = 15800
$ws.Range("K106").Value = 7504.385
$ws.Range("L106").Value = 15800
$ws.Range("M106").Value = -6873.385
$ws.Range("N106").Value = -17062

# Row 137 (Cutting Edge of Culinary Quality | Magnesia Whetstone)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 4404.448
$ws.Range("I137").Value = 3605.652
$ws.Range("J137").Value = 7466.5
$ws.Range("K137").Value = 10816.956
$ws.Range("L137").Value = 22399.5
$ws.Range("M137").Value = -8266.956
$ws.Range("N137").Value = -27499.5

# Row 138 (All-night Crafting | Cunning Craftsman's Tisane)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 7543.675
$ws.Range("J138").Value = 7386.0884
$ws.Range("L138").Value = 22158.2652
$ws.Range("N138").Value = -32438.2652

# Row 32 (Ingot We Trust | Steel Ingot)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3115.677
$ws.Range("I32").Value = 3055.8572
$ws.Range("K32").Value = 3055.8572
$ws.Range("M32").Value = -2768.8572

# Row 45 (Hollow Hallmarks | Mythril Ingot)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1767355.4
$ws.Range("I45").Value = 3403273.2
$ws.Range("K45").Value = 3403273.2
$ws.Range("M45").Value = -3402896.2

# Row 61 (Dealing with the Tough Stuff | Cobalt Ingot)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 8235930.5
$ws.Range("I61").Value = 11083027
$ws.Range("J61").Value = 2225393.2
$ws.Range("K61").Value = 11083027
$ws.Range("L61").Value = 2225393.2
$ws.Range("M61").Value = -11082815
$ws.Range("N61").Value = -2225817.2

# Row 110 (Scheduled Maintenance | Deepgold Ingot)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 6514.5454
$ws.Range("I110").Value = 7757.625
$ws.Range("J110").Value = 3199.6667
$ws.Range("K110").Value = 7757.625
$ws.Range("L110").Value = 3199.6667
$ws.Range("M110").Value = -5712.625
$ws.Range("N110").Value = -7289.6667

# Row 136 (Metal with Mettle | Cobalt Tungsten Ingot)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 8235930.5
$ws.Range("I136").Value = 11083027
$ws.Range("J136").Value = 2225393.2
$ws.Range("K136").Value = 33249081
$ws.Range("L136").Value = 6676179.600000001
$ws.Range("M136").Value = -33246531
$ws.Range("N136").Value = -6681279.600000001

# Row 134 (Ruthenium Supremium | Ruthenium Ingot)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4547605.5
$ws.Range("I134").Value = 1805.6316
$ws.Range("J134").Value = 33337672
$ws.Range("K134").Value = 5416.8948
$ws.Range("L134").Value = 100013016
$ws.Range("M134").Value = -2881.8948
$ws.Range("N134").Value = -100018086

# Row 31 (Wall Not Found | Walnut Lumber)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 16952482
$ws.Range("I31").Value = 28574640
$ws.Range("J31").Value = 3500
$ws.Range("K31").Value = 28574640
$ws.Range("L31").Value = 3500
$ws.Range("M31").Value = -28574345
$ws.Range("N31").Value = -4090

# Row 34 (Armoires of the Rich and Famous | Walnut Lumber)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 16952482
$ws.Range("I34").Value = 28574640
$ws.Range("J34").Value = 3500
$ws.Range("K34").Value = 28574640
$ws.Range("L34").Value = 3500
$ws.Range("M34").Value = -28574438
$ws.Range("N34").Value = -3904

# Row 93 (Reeling for Rods | Muudhorn Fishing Rod)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H93").Value = 53052.75
$ws.Range("I93").Value = 50603.062
$ws.Range("K93").Value = 50603.062
$ws.Range("M93").Value = -48731.062

# Row 107 (Built to Last | White Oak Lumber)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 913.36365
$ws.Range("I107").Value = 509.83334
$ws.Range("J107").Value = 2729.25
$ws.Range("K107").Value = 509.83334
$ws.Range("L107").Value = 2729.25
$ws.Range("M107").Value = 1410.16666
$ws.Range("N107").Value = -6569.25

# Row 5 (What a Sap | Maple Syrup)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1409.3
$ws.Range("I5").Value = 1165.8334
$ws.Range("J5").Value = 1774.5
$ws.Range("K5").Value = 3497.5002
$ws.Range("L5").Value = 5323.5
$ws.Range("M5").Value = -3385.5002
$ws.Range("N5").Value = -5547.5

# Row 29 (For Crumbs' Sake | Honey Muffin)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 6987.9287
$ws.Range("I29").Value = 2937.375
$ws.Range("J29").Value = 12388.667
$ws.Range("K29").Value = 8812.125
$ws.Range("L29").Value = 37166.001
$ws.Range("M29").Value = -8535.125
$ws.Range("N29").Value = -37720.001

# Row 68 (Such a Butter Face | Fermented Butter)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3020.1428
$ws.Range("J68").Value = 3235.75
$ws.Range("L68").Value = 9707.25
$ws.Range("N68").Value = -11329.25

# Row 71 (No Margarine of Error (L) | Fermented Butter)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 3020.1428
$ws.Range("J71").Value = 3235.75
$ws.Range("L71").Value = 29121.75
$ws.Range("N71").Value = -37233.75

# Row 104 (Fits to a Tea | Doman Tea)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 18888.834
$ws.Range("J104").Value = 18888.834
$ws.Range("L104").Value = 56666.50199999999
$ws.Range("N104").Value = -61908.50199999999

# Row 135 (Not-so-secret Ingredient | Royal Maple Syrup)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1409.3
$ws.Range("I135").Value = 1165.8334
$ws.Range("J135").Value = 1774.5
$ws.Range("K135").Value = 10492.5006
$ws.Range("L135").Value = 15970.5
$ws.Range("M135").Value = -7957.500599999999
$ws.Range("N135").Value = -21040.5

# Row 80 (Needs More Prayerbell | Hardsilver Ingot)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2386.0476
$ws.Range("I80").Value = 1677.7858
$ws.Range("J80").Value = 3802.5715
$ws.Range("K80").Value = 1677.7858
$ws.Range("L80").Value = 3802.5715
$ws.Range("M80").Value = -679.7858000000001
$ws.Range("N80").Value = -5798.5715

# Row 83 (With a Noise That Reaches Heaven (L) | Hardsilver Ingot)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2386.0476
$ws.Range("I83").Value = 1677.7858
$ws.Range("J83").Value = 3802.5715
$ws.Range("K83").Value = 8388.929
$ws.Range("L83").Value = 19012.8575
$ws.Range("M83").Value = -3396.929
$ws.Range("N83").Value = -28996.8575

# Row 61 (Spelling Me Softly | Raptor Leather)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 55558804
$ws.Range("I61").Value = 83336450
$ws.Range("J61").Value = 3511.8333
$ws.Range("K61").Value = 83336450
$ws.Range("L61").Value = 3511.8333
$ws.Range("M61").Value = -83336248
$ws.Range("N61").Value = -3915.8333

# Row 68 (You Could Say It's a Moving Target | Wyvern Leather)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2316964.5
$ws.Range("I68").Value = 3789494.2
$ws.Range("J68").Value = 2989.4285
$ws.Range("K68").Value = 3789494.2
$ws.Range("L68").Value = 2989.4285
$ws.Range("M68").Value = -3788745.2
$ws.Range("N68").Value = -4487.4285

# Row 71 (They Call It Bloody Mary (L) | Wyvern Leather)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2316964.5
$ws.Range("I71").Value = 3789494.2
$ws.Range("J71").Value = 2989.4285
$ws.Range("K71").Value = 18947471
$ws.Range("L71").Value = 14947.1425
$ws.Range("M71").Value = -18943727
$ws.Range("N71").Value = -22435.1425

# Row 82 (Trainin' the Neck | Dragon Leather)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 6448.7856
$ws.Range("I82").Value = 4998.2
$ws.Range("J82").Value = 7254.6665
$ws.Range("K82").Value = 4998.2
$ws.Range("L82").Value = 7254.6665
$ws.Range("M82").Value = -4637.2
$ws.Range("N82").Value = -7976.6665

# Row 85 (Training Is Only Skintight (L) | Dragon Leather)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 6448.7856
$ws.Range("I85").Value = 4998.2
$ws.Range("J85").Value = 7254.6665
$ws.Range("K85").Value = 4998.2
$ws.Range("L85").Value = 7254.6665
$ws.Range("M85").Value = -3750.2
$ws.Range("N85").Value = -9750.666499999999

# Row 113 (Peace in Rest | Atrociraptor Leather)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 55558804
$ws.Range("I113").Value = 83336450
$ws.Range("J113").Value = 3511.8333
$ws.Range("K113").Value = 83336450
$ws.Range("L113").Value = 3511.8333
$ws.Range("M113").Value = -83334280
$ws.Range("N113").Value = -7851.8333

# Row 122 (Hell on Leather | Gaja Leather)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2954.6326
$ws.Range("I122").Value = 2799.476
$ws.Range("J122").Value = 3885.5715
$ws.Range("K122").Value = 8398.428
$ws.Range("L122").Value = 11656.7145
$ws.Range("M122").Value = -5948.428
$ws.Range("N122").Value = -16556.7145

# Row 96 (Skills on Display | Ruby Cotton Cloth)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 7547
$ws.Range("I96").Value = 6330.727
$ws.Range("J96").Value = 9219.375
$ws.Range("K96").Value = 6330.727
$ws.Range("L96").Value = 9219.375
$ws.Range("M96").Value = -4957.727
$ws.Range("N96").Value = -11965.375
